$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(1).ColumnWidth = 12.333333333333334
$ws.Columns.Item(2).ColumnWidth = 10.166666666666666
$ws.Columns.Item(3).ColumnWidth = 22.666666666666668
$ws.Columns.Item(4).ColumnWidth = 12.333333333333334
$ws.Columns.Item(5).ColumnWidth = 36.833333333333336
$ws.Columns.Item(6).ColumnWidth = 32.666666666666664
$ws.Columns.Item(7).ColumnWidth = 28.833333333333332
$ws.Columns.Item(9).ColumnWidth = 12.666666666666666
